# Update gh-pages to output generated at 456a3b4
# This applies the "想去人数" (want-to-go count) refresh to the four
# worksheets of the workbook: 展览(1), 演出(2), 本地生活(3), 全部类型(4).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1211
$ws.Range("F5").Value = 1368
$ws.Range("F6").Value = 1713
$ws.Range("F7").Value = 6232
$ws.Range("F8").Value = 127
$ws.Range("F11").Value = 7
$ws.Range("F12").Value = 19
$ws.Range("F15").Value = 32
$ws.Range("F16").Value = 6971
$ws.Range("F21").Value = 1712
$ws.Range("F23").Value = 18
$ws.Range("F26").Value = 1597
$ws.Range("F27").Value = 766
$ws.Range("F28").Value = 318

# --- Sheet 2: 演出 (Performances) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F8").Value = 443
$ws.Range("F12").Value = 8
$ws.Range("F25").Value = 2

# --- Sheet 3: 本地生活 (Local Life) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F3").Value = 2260
$ws.Range("F4").Value = 663

# --- Sheet 4: 全部类型 (All Types) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 2260
$ws.Range("F4").Value = 663
$ws.Range("F5").Value = 1211
$ws.Range("F10").Value = 1368
$ws.Range("F12").Value = 1713
$ws.Range("F13").Value = 6232
$ws.Range("F19").Value = 19
$ws.Range("F21").Value = 8
$ws.Range("F23").Value = 6971
$ws.Range("F28").Value = 1712
$ws.Range("F30").Value = 18
$ws.Range("F33").Value = 1597
$ws.Range("F34").Value = 766
$ws.Range("F36").Value = 318
$ws.Range("F49").Value = 2
